# Implementing Student Report card Storing in Folder
#
# Refresh the automation login/test identifiers on the STAGE sheet:
#   - the School / Classroom / Section block used for the top ("automation"
#     / ditrictadmin) login row (A2:C2)
#   - the generated passwords for the fpk12admin / fpk12teacher /
#     fpk12student rows (E3:E5)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STAGE")

# --- A2:C2 -> new School / Classroom / Section names ------------------
$ws.Range("A2").Value = "FPK12School16006"
$ws.Range("B2").Value = "FPK12Classroom23482"
$ws.Range("C2").Value = "FPK12Section35176"

# Re-apply the centered / bordered look to the refreshed block.
$hdrRange = $ws.Range("A2:C2")
$hdrRange.HorizontalAlignment = -4108
$hdrRange.VerticalAlignment = -4108
$hdrRange.Borders.Item(9).LineStyle = 1
$hdrRange.Borders.Item(10).LineStyle = 1

# --- E3:E5 -> refreshed numeric passwords, stored as text -------------
# These columns hold numeric-looking ids as *text* (matching the rest of
# the column), so route the new value through a text-valued formula and
# convert it back to a literal via copy / paste-values instead of letting
# a bare numeric string auto-coerce to a number.
$e3 = $ws.Range("E3")
$e3.Formula = "=""40247"""
$e3.Copy()
$e3.PasteSpecial(-4163)

$e4 = $ws.Range("E4")
$e4.Formula = "=""67235"""
$e4.Copy()
$e4.PasteSpecial(-4163)

$e5 = $ws.Range("E5")
$e5.Formula = "=""60032"""
$e5.Copy()
$e5.PasteSpecial(-4163)

$wb.Application.CutCopyMode = 0
